# Update the "ASPE 2020 Annual Meeting" footer date from
#   "ASPE 2020 Annual Meeting – October 20-23, 2020"
# to
#   "ASPE 2020 Annual Meeting – October 21, 2020"
# (the last slide of the deck still carried the old wording; the other
# slides already read "October 21, 2020").

$p = $ppt.ActivePresentation

$oldText = "ASPE 2020 Annual Meeting " + [char]0x2013 + " October 20-23, 2020"
$newText = "ASPE 2020 Annual Meeting " + [char]0x2013 + " October 21, 2020"

for ($s = 1; $s -le $p.Slides.Count; $s++) {
    $slide = $p.Slides.Item($s)
    for ($i = 1; $i -le $slide.Shapes.Count; $i++) {
        $shape = $slide.Shapes.Item($i)
        if ($shape.HasTextFrame) {
            $tr = $shape.TextFrame.TextRange
            if ($tr.Text -eq $oldText) {
                $tr.Text = $newText
            }
        }
    }
}
